$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 31, shifting existing rows 31..112 down to 32..113
$ws.Rows.Item(31).Insert()

# Populate the new row 31 with the new weekly data point.
# Columns A, B, C, E, F, G, H, I, N, O, Q, R are constant across all rows in
# this sheet, so copy them from the row above (row 30) to keep formatting
# and content consistent; only D, J, K, L, M, P hold the new record's data.
$ws.Range("A31:R31").Value = $ws.Range("A30:R30").Value2

$ws.Range("D31").Value = 44459
$ws.Range("J31").Value = 2800
$ws.Range("K31").Value = 900
$ws.Range("L31").Value = 1000
$ws.Range("M31").Value = 950
$ws.Range("P31").Value = 158
